# This sheet contains weekly "Alcachofa" price records. A new, more recent
# weekly record is inserted at row 41 (pushing the existing rows 41-122 down
# to rows 42-123), matching the source data feed's weekly update pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; this shifts rows 41:122 down to 42:123
# and extends the used range to A1:R123.
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new weekly record.
$ws.Cells.Item(41, 1).Value  = 7
$ws.Cells.Item(41, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(41, 3).Value  = "Ñuble"
$ws.Cells.Item(41, 4).Value  = 45162
$ws.Cells.Item(41, 5).Value  = 16
$ws.Cells.Item(41, 6).Value  = 100112013
$ws.Cells.Item(41, 7).Value  = "Alcachofa"
$ws.Cells.Item(41, 8).Value  = "Madrigal"
$ws.Cells.Item(41, 9).Value  = "Primera"
$ws.Cells.Item(41, 10).Value = 300
$ws.Cells.Item(41, 11).Value = 11000
$ws.Cells.Item(41, 12).Value = 11000
$ws.Cells.Item(41, 13).Value = 11000
$ws.Cells.Item(41, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(41, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(41, 16).Value = 275
$ws.Cells.Item(41, 17).Value = 40
$ws.Cells.Item(41, 18).Value = "Hortaliza"
